$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the manager/password test data so every row (TC_Valid1..3) now shares
# the same refreshed credentials (mngr122407 / apYsUtA), replacing the old
# per-row mngr11914x / random-password values.
$ws.Range("B2").Value = "mngr122407"
$ws.Range("C2").Value = "apYsUtA"
$ws.Range("B3").Value = "mngr122407"
$ws.Range("C3").Value = "apYsUtA"
$ws.Range("B4").Value = "mngr122407"
$ws.Range("C4").Value = "apYsUtA"

# Widen columns B (UserName) and C (Password) to fit the new values.
$ws.Columns("B").ColumnWidth = 24
$ws.Columns("C").ColumnWidth = 17.333333333333332

# Move the active selection off the data range, as captured in the saved view.
$ws.Range("E5").Select()
